# Add I0 and IF columns to the sheet (columns I and J), mirroring the
# header style used by the existing headers (B1:H1), and fill in the
# corresponding data values for rows 2-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell (H1) onto the new headers
# so they match the bold/centered/bordered look of the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for I2:J28 (I0, IF) taken from the diff
$data = @(
    @(8, 8),
    @(5, 5),
    @(7, 8),
    @(7, 8),
    @(1, 3),
    @(1, 4),
    @(7, 8),
    @(3, 7),
    @(2, 6),
    @(10, 10),
    @(6, 6),
    @(3, 5),
    @(7, 9),
    @(5, 6),
    @(5, 7),
    @(5, 6),
    @(1, 3),
    @(1, 4),
    @(1, 5),
    @(1, 5),
    @(7, 7),
    @(6, 7),
    @(1, 6),
    @(4, 8),
    @(6, 8),
    @(4, 5),
    @(3, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]  # Column I
    $ws.Cells.Item($row, 10).Value = $pair[1] # Column J
    $row++
}
